# "Generate Report for Archive"
#
# The localization pipeline re-ran: the handoff files moved from
# "Ready for handoff" into "In Translation", and (because that status text
# is shorter) the Status columns got narrower on every sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E and F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
